$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header typo in B1: "montly_sales" -> "monthly_sales"
# (A1 "advertising_expense" is left untouched)
$ws.Range("B1").Value = "monthly_sales"

# Overwrite a batch of previously-formula-driven "monthly sales" cells in
# column B with hand-entered numbers (rounded estimates replacing the
# computed A*2.5 figures). The remaining B cells keep their original
# "=A#*2.5" formula.
$ws.Range("B3").Value = 7000
$ws.Range("B5").Value = 3500
$ws.Range("B6").Value = 4100
$ws.Range("B8").Value = 4050
$ws.Range("B11").Value = 7400
$ws.Range("B12").Value = 6050
$ws.Range("B13").Value = 3200
$ws.Range("B14").Value = 7227
$ws.Range("B15").Value = 4669
$ws.Range("B16").Value = 4402
$ws.Range("B17").Value = 6000
$ws.Range("B19").Value = 6030
$ws.Range("B21").Value = 6550
$ws.Range("B22").Value = 6920
$ws.Range("B24").Value = 6200
$ws.Range("B25").Value = 2700
$ws.Range("B26").Value = 6600
$ws.Range("B27").Value = 4100
$ws.Range("B28").Value = 6500
$ws.Range("B29").Value = 7300

# Move the active selection, matching the author's saved cursor position.
$ws.Range("P13").Select() | Out-Null
